# Recipe_Search_Questions.xlsx — add two new recipe-search question rows
# (rows 43 and 44) to Sheet1, matching the "ollama is working !!!" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: Semantica / Pergunta difícil + / new query
$ws.Range("A43").Value = "Semantica"
$ws.Range("B43").Value = "Pergunta difícil +"
$ws.Range("C43").Value = "I have 9min maximum to make a lunch, can you help me?"

# Row 44: Semantica / Pergunta difícil + / new query
$ws.Range("A44").Value = "Semantica"
$ws.Range("B44").Value = "Pergunta difícil +"
$ws.Range("C44").Value = "Esfiha de carne vegana"

# Reflect the author's final cursor position (selection ends on B44)
$null = $ws.Range("B44").Select()
